$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("End point")

# Update max volume column B for rows 17-23, and clear C/D columns 16-23
$ws.Range("B17").Value = 20000
$ws.Range("B18").Value = 10000
$ws.Range("B19").Value = 5000
$ws.Range("B20").Value = 4000
$ws.Range("B21").Value = 3000
$ws.Range("B22").Value = 2000
$ws.Range("B23").Value = 1000

$ws.Range("C16:D23").ClearContents()

$ws.Range("C26").Select()
